# Auto-generated Excel COM-interop script to apply the "Horarios actualizados" update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# Header rows (timestamp + row-count banner)
$ws.Range("A2").Value = "Última actualización: 16:52:36"
$ws.Range("A3").Value = "Total filas: 359"
$ws2.Range("A2").Value = "Última actualización: 16:52:36"
$ws3.Range("A2").Value = "Última actualización: 16:52:36"

# Refreshed schedule data, rows 6-364 (A:Hora_Scrap B:Hora_Llegada C:Linea D:Minutos E:Parada)
$data = New-Object 'object[,]' 359,5
$data[0,0] = "03:45:25"
$data[0,1] = "03:47"
$data[0,2] = "14_ABASTO"
$data[0,3] = 2
$data[0,4] = "LP1912"
$data[1,0] = "03:45:25"
$data[1,1] = "04:01"
$data[1,2] = "81_EL PELIGRO"
$data[1,3] = 16
$data[1,4] = "LP1912"
$data[2,0] = "03:45:25"
$data[2,1] = "04:46"
$data[2,2] = "215A_EL PATO"
$data[2,3] = 61
$data[2,4] = "LP1912"
$data[3,0] = "03:45:25"
$data[3,1] = "04:53"
$data[3,2] = "11_ETCHEVERRY"
$data[3,3] = 68
$data[3,4] = "LP1912"
$data[4,0] = "04:56:49"
$data[4,1] = "05:13"
$data[4,2] = "14_ABASTO"
$data[4,3] = 17
$data[4,4] = "LP1912"
$data[5,0] = "03:45:25"
$data[5,1] = "05:16"
$data[5,2] = "17_ROMERO"
$data[5,3] = 91
$data[5,4] = "LP1912"
$data[6,0] = "04:45:05"
$data[6,1] = "05:16"
$data[6,2] = "14_ABASTO"
$data[6,3] = 31
$data[6,4] = "LP1912"
$data[7,0] = "03:45:25"
$data[7,1] = "05:22"
$data[7,2] = "23_HERNANDEZ"
$data[7,3] = 97
$data[7,4] = "LP1912"
$data[8,0] = "05:26:08"
$data[8,1] = "05:28"
$data[8,2] = "14_ABASTO"
$data[8,3] = 2
$data[8,4] = "LP1912"
$data[9,0] = "04:18:02"
$data[9,1] = "05:34"
$data[9,2] = "14_ABASTO"
$data[9,3] = 76
$data[9,4] = "LP1912"
$data[10,0] = "03:45:25"
$data[10,1] = "05:34"
$data[10,2] = "215B_EL PATO"
$data[10,3] = 109
$data[10,4] = "LP1912"
$data[11,0] = "04:18:02"
$data[11,1] = "05:35"
$data[11,2] = "215B_EL PATO"
$data[11,3] = 77
$data[11,4] = "LP1912"
$data[12,0] = "03:45:25"
$data[12,1] = "05:37"
$data[12,2] = "14_ABASTO"
$data[12,3] = 112
$data[12,4] = "LP1912"
$data[13,0] = "04:18:02"
$data[13,1] = "05:46"
$data[13,2] = "15_ABASTO"
$data[13,3] = 88
$data[13,4] = "LP1912"
$data[14,0] = "04:45:05"
$data[14,1] = "06:04"
$data[14,2] = "16_SANTA ANA"
$data[14,3] = 79
$data[14,4] = "LP1912"
$data[15,0] = "04:18:02"
$data[15,1] = "06:05"
$data[15,2] = "16_SANTA ANA"
$data[15,3] = 107
$data[15,4] = "LP1912"
$data[16,0] = "04:56:49"
$data[16,1] = "06:11"
$data[16,2] = "215A_EL PATO"
$data[16,3] = 75
$data[16,4] = "LP1912"
$data[17,0] = "04:18:02"
$data[17,1] = "06:12"
$data[17,2] = "215A_EL PATO"
$data[17,3] = 114
$data[17,4] = "LP1912"
$data[18,0] = "04:18:02"
$data[18,1] = "06:14"
$data[18,2] = "225_HARAS DEL SUR"
$data[18,3] = 116
$data[18,4] = "LP1912"
$data[19,0] = "04:45:05"
$data[19,1] = "06:21"
$data[19,2] = "26_HERNANDEZ"
$data[19,3] = 96
$data[19,4] = "LP1912"
$data[20,0] = "06:25:43"
$data[20,1] = "06:26"
$data[20,2] = "86_EST CHICA-ESC AGRARIA"
$data[20,3] = 1
$data[20,4] = "LP1912"
$data[21,0] = "04:45:05"
$data[21,1] = "06:27"
$data[21,2] = "23_HERNANDEZ"
$data[21,3] = 102
$data[21,4] = "LP1912"
$data[22,0] = "06:25:43"
$data[22,1] = "06:28"
$data[22,2] = "23_HERNANDEZ"
$data[22,3] = 3
$data[22,4] = "LP1912"
$data[23,0] = "04:56:49"
$data[23,1] = "06:29"
$data[23,2] = "86_EST CHICA-ESC AGRARIA"
$data[23,3] = 93
$data[23,4] = "LP1912"
$data[24,0] = "04:45:05"
$data[24,1] = "06:30"
$data[24,2] = "86_EST CHICA-ESC AGRARIA"
$data[24,3] = 105
$data[24,4] = "LP1912"
$data[25,0] = "04:45:05"
$data[25,1] = "06:31"
$data[25,2] = "16_SANTA ANA"
$data[25,3] = 106
$data[25,4] = "LP1912"
$data[26,0] = "05:55:25"
$data[26,1] = "06:44"
$data[26,2] = "26_HERNANDEZ"
$data[26,3] = 49
$data[26,4] = "LP1912"
$data[27,0] = "04:45:05"
$data[27,1] = "06:44"
$data[27,2] = "225_C ROCA-H SUR"
$data[27,3] = 119
$data[27,4] = "LP1912"
$data[28,0] = "04:56:49"
$data[28,1] = "06:46"
$data[28,2] = "215C_EL PATO"
$data[28,3] = 110
$data[28,4] = "LP1912"
$data[29,0] = "05:26:08"
$data[29,1] = "06:47"
$data[29,2] = "215C_EL PATO"
$data[29,3] = 81
$data[29,4] = "LP1912"
$data[30,0] = "05:55:25"
$data[30,1] = "06:59"
$data[30,2] = "14_ABASTO"
$data[30,3] = 64
$data[30,4] = "LP1912"
$data[31,0] = "05:26:08"
$data[31,1] = "07:00"
$data[31,2] = "14_ABASTO"
$data[31,3] = 94
$data[31,4] = "LP1912"
$data[32,0] = "06:25:43"
$data[32,1] = "07:01"
$data[32,2] = "16_SANTA ANA"
$data[32,3] = 36
$data[32,4] = "LP1912"
$data[33,0] = "05:55:25"
$data[33,1] = "07:04"
$data[33,2] = "23_HERNANDEZ"
$data[33,3] = 69
$data[33,4] = "LP1912"
$data[34,0] = "05:26:08"
$data[34,1] = "07:05"
$data[34,2] = "23_HERNANDEZ"
$data[34,3] = 99
$data[34,4] = "LP1912"
$data[35,0] = "05:26:08"
$data[35,1] = "07:05"
$data[35,2] = "15_ABASTO"
$data[35,3] = 99
$data[35,4] = "LP1912"
$data[36,0] = "05:26:08"
$data[36,1] = "07:06"
$data[36,2] = "10_OLMOS"
$data[36,3] = 100
$data[36,4] = "LP1912"
$data[37,0] = "05:26:08"
$data[37,1] = "07:07"
$data[37,2] = "225_GOMEZ"
$data[37,3] = 101
$data[37,4] = "LP1912"
$data[38,0] = "05:26:08"
$data[38,1] = "07:11"
$data[38,2] = "215A_EL PATO"
$data[38,3] = 105
$data[38,4] = "LP1912"
$data[39,0] = "06:55:02"
$data[39,1] = "07:12"
$data[39,2] = "215A_EL PATO"
$data[39,3] = 17
$data[39,4] = "LP1912"
$data[40,0] = "06:25:43"
$data[40,1] = "07:14"
$data[40,2] = "26_HERNANDEZ"
$data[40,3] = 49
$data[40,4] = "LP1912"
$data[41,0] = "05:55:25"
$data[41,1] = "07:15"
$data[41,2] = "11_ETCHEVERRY"
$data[41,3] = 80
$data[41,4] = "LP1912"
$data[42,0] = "05:26:08"
$data[42,1] = "07:16"
$data[42,2] = "11_ETCHEVERRY"
$data[42,3] = 110
$data[42,4] = "LP1912"
$data[43,0] = "06:55:02"
$data[43,1] = "07:17"
$data[43,2] = "16_SANTA ANA"
$data[43,3] = 22
$data[43,4] = "LP1912"
$data[44,0] = "05:26:08"
$data[44,1] = "07:21"
$data[44,2] = "26_HERNANDEZ"
$data[44,3] = 115
$data[44,4] = "LP1912"
$data[45,0] = "05:26:08"
$data[45,1] = "07:23"
$data[45,2] = "10_OLMOS"
$data[45,3] = 117
$data[45,4] = "LP1912"
$data[46,0] = "05:55:25"
$data[46,1] = "07:30"
$data[46,2] = "10_OLMOS"
$data[46,3] = 95
$data[46,4] = "LP1912"
$data[47,0] = "05:55:25"
$data[47,1] = "07:31"
$data[47,2] = "11_ETCHEVERRY"
$data[47,3] = 96
$data[47,4] = "LP1912"
$data[48,0] = "05:55:25"
$data[48,1] = "07:31"
$data[48,2] = "16_SANTA ANA"
$data[48,3] = 96
$data[48,4] = "LP1912"
$data[49,0] = "06:55:02"
$data[49,1] = "07:32"
$data[49,2] = "11_ETCHEVERRY"
$data[49,3] = 37
$data[49,4] = "LP1912"
$data[50,0] = "05:55:25"
$data[50,1] = "07:32"
$data[50,2] = "84_COLONIA URQUIZA-ESC 49"
$data[50,3] = 97
$data[50,4] = "LP1912"
$data[51,0] = "06:55:02"
$data[51,1] = "07:32"
$data[51,2] = "16_SANTA ANA"
$data[51,3] = 37
$data[51,4] = "LP1912"
$data[52,0] = "07:19:29"
$data[52,1] = "07:35"
$data[52,2] = "23_HERNANDEZ"
$data[52,3] = 16
$data[52,4] = "LP1912"
$data[53,0] = "05:55:25"
$data[53,1] = "07:36"
$data[53,2] = "27_EL RETIRO"
$data[53,3] = 101
$data[53,4] = "LP1912"
$data[54,0] = "06:55:02"
$data[54,1] = "07:37"
$data[54,2] = "27_EL RETIRO"
$data[54,3] = 42
$data[54,4] = "LP1912"
$data[55,0] = "05:55:25"
$data[55,1] = "07:39"
$data[55,2] = "10_OLMOS"
$data[55,3] = 104
$data[55,4] = "LP1912"
$data[56,0] = "07:19:29"
$data[56,1] = "07:46"
$data[56,2] = "16_SANTA ANA"
$data[56,3] = 27
$data[56,4] = "LP1912"
$data[57,0] = "05:55:25"
$data[57,1] = "07:47"
$data[57,2] = "14_ABASTO"
$data[57,3] = 112
$data[57,4] = "LP1912"
$data[58,0] = "06:55:02"
$data[58,1] = "07:48"
$data[58,2] = "14_ABASTO"
$data[58,3] = 53
$data[58,4] = "LP1912"
$data[59,0] = "07:50:16"
$data[59,1] = "07:50"
$data[59,2] = "10_OLMOS"
$data[59,3] = 0
$data[59,4] = "LP1912"
$data[60,0] = "05:55:25"
$data[60,1] = "07:51"
$data[60,2] = "215D_EL PATO"
$data[60,3] = 116
$data[60,4] = "LP1912"
$data[61,0] = "06:55:02"
$data[61,1] = "07:52"
$data[61,2] = "215D_EL PATO"
$data[61,3] = 57
$data[61,4] = "LP1912"
$data[62,0] = "07:19:29"
$data[62,1] = "07:59"
$data[62,2] = "23_HERNANDEZ"
$data[62,3] = 40
$data[62,4] = "LP1912"
$data[63,0] = "06:25:43"
$data[63,1] = "08:01"
$data[63,2] = "23_HERNANDEZ"
$data[63,3] = 96
$data[63,4] = "LP1912"
$data[64,0] = "07:19:29"
$data[64,1] = "08:03"
$data[64,2] = "11_ETCHEVERRY"
$data[64,3] = 44
$data[64,4] = "LP1912"
$data[65,0] = "06:55:02"
$data[65,1] = "08:03"
$data[65,2] = "23_HERNANDEZ"
$data[65,3] = 68
$data[65,4] = "LP1912"
$data[66,0] = "08:02:22"
$data[66,1] = "08:05"
$data[66,2] = "23_HERNANDEZ"
$data[66,3] = 3
$data[66,4] = "LP1912"
$data[67,0] = "07:19:29"
$data[67,1] = "08:10"
$data[67,2] = "16_SANTA ANA"
$data[67,3] = 51
$data[67,4] = "LP1912"
$data[68,0] = "07:50:16"
$data[68,1] = "08:11"
$data[68,2] = "16_SANTA ANA"
$data[68,3] = 21
$data[68,4] = "LP1912"
$data[69,0] = "06:25:43"
$data[69,1] = "08:12"
$data[69,2] = "15_ABASTO"
$data[69,3] = 107
$data[69,4] = "LP1912"
$data[70,0] = "07:50:16"
$data[70,1] = "08:13"
$data[70,2] = "10_OLMOS"
$data[70,3] = 23
$data[70,4] = "LP1912"
$data[71,0] = "06:55:02"
$data[71,1] = "08:21"
$data[71,2] = "26_HERNANDEZ"
$data[71,3] = 86
$data[71,4] = "LP1912"
$data[72,0] = "06:25:43"
$data[72,1] = "08:22"
$data[72,2] = "16_P MOR-SANTA ANA"
$data[72,3] = 117
$data[72,4] = "LP1912"
$data[73,0] = "06:55:02"
$data[73,1] = "08:23"
$data[73,2] = "16_P MOR-SANTA ANA"
$data[73,3] = 88
$data[73,4] = "LP1912"
$data[74,0] = "06:25:43"
$data[74,1] = "08:23"
$data[74,2] = "215B_EL PATO"
$data[74,3] = 118
$data[74,4] = "LP1912"
$data[75,0] = "06:55:02"
$data[75,1] = "08:27"
$data[75,2] = "84_COLONIA URQUIZA-ESC 49"
$data[75,3] = 92
$data[75,4] = "LP1912"
$data[76,0] = "07:50:16"
$data[76,1] = "08:30"
$data[76,2] = "23_HERNANDEZ"
$data[76,3] = 40
$data[76,4] = "LP1912"
$data[77,0] = "08:02:22"
$data[77,1] = "08:33"
$data[77,2] = "10_OLMOS"
$data[77,3] = 31
$data[77,4] = "LP1912"
$data[78,0] = "08:02:22"
$data[78,1] = "08:34"
$data[78,2] = "23_HERNANDEZ"
$data[78,3] = 32
$data[78,4] = "LP1912"
$data[79,0] = "08:32:09"
$data[79,1] = "08:37"
$data[79,2] = "23_HERNANDEZ"
$data[79,3] = 5
$data[79,4] = "LP1912"
$data[80,0] = "06:55:02"
$data[80,1] = "08:42"
$data[80,2] = "81_EL PELIGRO"
$data[80,3] = 107
$data[80,4] = "LP1912"
$data[81,0] = "07:19:29"
$data[81,1] = "08:43"
$data[81,2] = "14_ABASTO"
$data[81,3] = 84
$data[81,4] = "LP1912"
$data[82,0] = "07:50:16"
$data[82,1] = "08:44"
$data[82,2] = "14_ABASTO"
$data[82,3] = 54
$data[82,4] = "LP1912"
$data[83,0] = "08:32:09"
$data[83,1] = "08:53"
$data[83,2] = "10_OLMOS"
$data[83,3] = 21
$data[83,4] = "LP1912"
$data[84,0] = "06:55:02"
$data[84,1] = "08:54"
$data[84,2] = "17_ROMERO"
$data[84,3] = 119
$data[84,4] = "LP1912"
$data[85,0] = "07:19:29"
$data[85,1] = "09:01"
$data[85,2] = "215A_EL PATO"
$data[85,3] = 102
$data[85,4] = "LP1912"
$data[86,0] = "07:50:16"
$data[86,1] = "09:02"
$data[86,2] = "215A_EL PATO"
$data[86,3] = 72
$data[86,4] = "LP1912"
$data[87,0] = "08:02:22"
$data[87,1] = "09:03"
$data[87,2] = "11_ETCHEVERRY"
$data[87,3] = 61
$data[87,4] = "LP1912"
$data[88,0] = "08:32:09"
$data[88,1] = "09:04"
$data[88,2] = "11_ETCHEVERRY"
$data[88,3] = 32
$data[88,4] = "LP1912"
$data[89,0] = "08:32:09"
$data[89,1] = "09:05"
$data[89,2] = "23_HERNANDEZ"
$data[89,3] = 33
$data[89,4] = "LP1912"
$data[90,0] = "07:19:29"
$data[90,1] = "09:10"
$data[90,2] = "16_P MOR-SANTA ANA"
$data[90,3] = 111
$data[90,4] = "LP1912"
$data[91,0] = "07:50:16"
$data[91,1] = "09:11"
$data[91,2] = "16_P MOR-SANTA ANA"
$data[91,3] = 81
$data[91,4] = "LP1912"
$data[92,0] = "08:32:09"
$data[92,1] = "09:13"
$data[92,2] = "10_OLMOS"
$data[92,3] = 41
$data[92,4] = "LP1912"
$data[93,0] = "07:19:29"
$data[93,1] = "09:16"
$data[93,2] = "27_EL RETIRO"
$data[93,3] = 117
$data[93,4] = "LP1912"
$data[94,0] = "07:50:16"
$data[94,1] = "09:17"
$data[94,2] = "27_EL RETIRO"
$data[94,3] = 87
$data[94,4] = "LP1912"
$data[95,0] = "07:50:16"
$data[95,1] = "09:21"
$data[95,2] = "26_HERNANDEZ"
$data[95,3] = 91
$data[95,4] = "LP1912"
$data[96,0] = "08:02:22"
$data[96,1] = "09:22"
$data[96,2] = "16_SANTA ANA"
$data[96,3] = 80
$data[96,4] = "LP1912"
$data[97,0] = "08:32:09"
$data[97,1] = "09:23"
$data[97,2] = "16_SANTA ANA"
$data[97,3] = 51
$data[97,4] = "LP1912"
$data[98,0] = "07:50:16"
$data[98,1] = "09:23"
$data[98,2] = "17_ROMERO"
$data[98,3] = 93
$data[98,4] = "LP1912"
$data[99,0] = "08:02:22"
$data[99,1] = "09:23"
$data[99,2] = "11_ETCHEVERRY"
$data[99,3] = 81
$data[99,4] = "LP1912"
$data[100,0] = "07:50:16"
$data[100,1] = "09:24"
$data[100,2] = "11_ETCHEVERRY"
$data[100,3] = 94
$data[100,4] = "LP1912"
$data[101,0] = "07:50:16"
$data[101,1] = "09:28"
$data[101,2] = "16_SANTA ANA"
$data[101,3] = 98
$data[101,4] = "LP1912"
$data[102,0] = "07:50:16"
$data[102,1] = "09:32"
$data[102,2] = "15_ABASTO"
$data[102,3] = 102
$data[102,4] = "LP1912"
$data[103,0] = "07:50:16"
$data[103,1] = "09:33"
$data[103,2] = "10_OLMOS"
$data[103,3] = 103
$data[103,4] = "LP1912"
$data[104,0] = "08:56:29"
$data[104,1] = "09:34"
$data[104,2] = "23_HERNANDEZ"
$data[104,3] = 38
$data[104,4] = "LP1912"
$data[105,0] = "08:56:29"
$data[105,1] = "09:34"
$data[105,2] = "16_SANTA ANA"
$data[105,3] = 38
$data[105,4] = "LP1912"
$data[106,0] = "08:32:09"
$data[106,1] = "09:35"
$data[106,2] = "16_SANTA ANA"
$data[106,3] = 63
$data[106,4] = "LP1912"
$data[107,0] = "08:48:08"
$data[107,1] = "09:35"
$data[107,2] = "23_HERNANDEZ"
$data[107,3] = 47
$data[107,4] = "LP1912"
$data[108,0] = "09:35:26"
$data[108,1] = "09:39"
$data[108,2] = "23_HERNANDEZ"
$data[108,3] = 4
$data[108,4] = "LP1912"
$data[109,0] = "07:50:16"
$data[109,1] = "09:42"
$data[109,2] = "215C_EL PATO"
$data[109,3] = 112
$data[109,4] = "LP1912"
$data[110,0] = "08:02:22"
$data[110,1] = "09:43"
$data[110,2] = "14_ABASTO"
$data[110,3] = 101
$data[110,4] = "LP1912"
$data[111,0] = "07:50:16"
$data[111,1] = "09:44"
$data[111,2] = "14_ABASTO"
$data[111,3] = 114
$data[111,4] = "LP1912"
$data[112,0] = "09:35:26"
$data[112,1] = "09:46"
$data[112,2] = "16_SANTA ANA"
$data[112,3] = 11
$data[112,4] = "LP1912"
$data[113,0] = "08:32:09"
$data[113,1] = "09:52"
$data[113,2] = "15_ABASTO"
$data[113,3] = 80
$data[113,4] = "LP1912"
$data[114,0] = "08:56:29"
$data[114,1] = "09:53"
$data[114,2] = "10_OLMOS"
$data[114,3] = 57
$data[114,4] = "LP1912"
$data[115,0] = "09:35:26"
$data[115,1] = "09:58"
$data[115,2] = "16_SANTA ANA"
$data[115,3] = 23
$data[115,4] = "LP1912"
$data[116,0] = "09:35:26"
$data[116,1] = "10:03"
$data[116,2] = "11_ETCHEVERRY"
$data[116,3] = 28
$data[116,4] = "LP1912"
$data[117,0] = "08:56:29"
$data[117,1] = "10:10"
$data[117,2] = "16_P MOR-SANTA ANA"
$data[117,3] = 74
$data[117,4] = "LP1912"
$data[118,0] = "08:32:09"
$data[118,1] = "10:11"
$data[118,2] = "16_P MOR-SANTA ANA"
$data[118,3] = 99
$data[118,4] = "LP1912"
$data[119,0] = "09:35:26"
$data[119,1] = "10:12"
$data[119,2] = "15_ABASTO"
$data[119,3] = 37
$data[119,4] = "LP1912"
$data[120,0] = "09:35:26"
$data[120,1] = "10:13"
$data[120,2] = "10_OLMOS"
$data[120,3] = 38
$data[120,4] = "LP1912"
$data[121,0] = "08:32:09"
$data[121,1] = "10:21"
$data[121,2] = "26_HERNANDEZ"
$data[121,3] = 109
$data[121,4] = "LP1912"
$data[122,0] = "08:32:09"
$data[122,1] = "10:22"
$data[122,2] = "17_ROMERO"
$data[122,3] = 110
$data[122,4] = "LP1912"
$data[123,0] = "09:35:26"
$data[123,1] = "10:23"
$data[123,2] = "11_ETCHEVERRY"
$data[123,3] = 48
$data[123,4] = "LP1912"
$data[124,0] = "08:56:29"
$data[124,1] = "10:26"
$data[124,2] = "215A_EL PATO"
$data[124,3] = 90
$data[124,4] = "LP1912"
$data[125,0] = "08:32:09"
$data[125,1] = "10:27"
$data[125,2] = "215A_EL PATO"
$data[125,3] = 115
$data[125,4] = "LP1912"
$data[126,0] = "10:29:57"
$data[126,1] = "10:29"
$data[126,2] = "16_SANTA ANA"
$data[126,3] = 0
$data[126,4] = "LP1912"
$data[127,0] = "10:29:57"
$data[127,1] = "10:31"
$data[127,2] = "10_OLMOS"
$data[127,3] = 2
$data[127,4] = "LP1912"
$data[128,0] = "09:35:26"
$data[128,1] = "10:34"
$data[128,2] = "23_HERNANDEZ"
$data[128,3] = 59
$data[128,4] = "LP1912"
$data[129,0] = "10:29:57"
$data[129,1] = "10:34"
$data[129,2] = "16_SANTA ANA"
$data[129,3] = 5
$data[129,4] = "LP1912"
$data[130,0] = "10:29:57"
$data[130,1] = "10:39"
$data[130,2] = "23_HERNANDEZ"
$data[130,3] = 10
$data[130,4] = "LP1912"
$data[131,0] = "10:29:57"
$data[131,1] = "10:41"
$data[131,2] = "17_ROMERO"
$data[131,3] = 12
$data[131,4] = "LP1912"
$data[132,0] = "08:48:08"
$data[132,1] = "10:42"
$data[132,2] = "17_ROMERO"
$data[132,3] = 114
$data[132,4] = "LP1912"
$data[133,0] = "08:56:29"
$data[133,1] = "10:43"
$data[133,2] = "14_ABASTO"
$data[133,3] = 107
$data[133,4] = "LP1912"
$data[134,0] = "08:48:08"
$data[134,1] = "10:44"
$data[134,2] = "14_ABASTO"
$data[134,3] = 116
$data[134,4] = "LP1912"
$data[135,0] = "10:29:57"
$data[135,1] = "10:51"
$data[135,2] = "15_ABASTO"
$data[135,3] = 22
$data[135,4] = "LP1912"
$data[136,0] = "10:29:57"
$data[136,1] = "10:52"
$data[136,2] = "10_OLMOS"
$data[136,3] = 23
$data[136,4] = "LP1912"
$data[137,0] = "09:35:26"
$data[137,1] = "10:54"
$data[137,2] = "27_EL RETIRO"
$data[137,3] = 79
$data[137,4] = "LP1912"
$data[138,0] = "10:29:57"
$data[138,1] = "10:56"
$data[138,2] = "27_EL RETIRO"
$data[138,3] = 27
$data[138,4] = "LP1912"
$data[139,0] = "10:59:49"
$data[139,1] = "10:59"
$data[139,2] = "16_SANTA ANA"
$data[139,3] = 0
$data[139,4] = "LP1912"
$data[140,0] = "10:29:57"
$data[140,1] = "11:01"
$data[140,2] = "215C_EL PATO"
$data[140,3] = 32
$data[140,4] = "LP1912"
$data[141,0] = "09:35:26"
$data[141,1] = "11:02"
$data[141,2] = "215C_EL PATO"
$data[141,3] = 87
$data[141,4] = "LP1912"
$data[142,0] = "10:29:57"
$data[142,1] = "11:03"
$data[142,2] = "11_ETCHEVERRY"
$data[142,3] = 34
$data[142,4] = "LP1912"
$data[143,0] = "10:29:57"
$data[143,1] = "11:04"
$data[143,2] = "23_HERNANDEZ"
$data[143,3] = 35
$data[143,4] = "LP1912"
$data[144,0] = "10:59:49"
$data[144,1] = "11:06"
$data[144,2] = "23_HERNANDEZ"
$data[144,3] = 7
$data[144,4] = "LP1912"
$data[145,0] = "09:35:26"
$data[145,1] = "11:06"
$data[145,2] = "16_P MOR-167 Y 521"
$data[145,3] = 91
$data[145,4] = "LP1912"
$data[146,0] = "10:59:49"
$data[146,1] = "11:11"
$data[146,2] = "10_OLMOS"
$data[146,3] = 12
$data[146,4] = "LP1912"
$data[147,0] = "10:29:57"
$data[147,1] = "11:11"
$data[147,2] = "15_ABASTO"
$data[147,3] = 42
$data[147,4] = "LP1912"
$data[148,0] = "10:59:49"
$data[148,1] = "11:12"
$data[148,2] = "15_ABASTO"
$data[148,3] = 13
$data[148,4] = "LP1912"
$data[149,0] = "09:35:26"
$data[149,1] = "11:19"
$data[149,2] = "86_EST CHICA-ESC AGRARIA"
$data[149,3] = 104
$data[149,4] = "LP1912"
$data[150,0] = "10:29:57"
$data[150,1] = "11:20"
$data[150,2] = "26_HERNANDEZ"
$data[150,3] = 51
$data[150,4] = "LP1912"
$data[151,0] = "09:35:26"
$data[151,1] = "11:21"
$data[151,2] = "26_HERNANDEZ"
$data[151,3] = 106
$data[151,4] = "LP1912"
$data[152,0] = "10:29:57"
$data[152,1] = "11:26"
$data[152,2] = "225_C ROCA-H SUR"
$data[152,3] = 57
$data[152,4] = "LP1912"
$data[153,0] = "09:35:26"
$data[153,1] = "11:27"
$data[153,2] = "225_C ROCA-H SUR"
$data[153,3] = 112
$data[153,4] = "LP1912"
$data[154,0] = "11:30:45"
$data[154,1] = "11:30"
$data[154,2] = "16_SANTA ANA"
$data[154,3] = 0
$data[154,4] = "LP1912"
$data[155,0] = "11:30:45"
$data[155,1] = "11:31"
$data[155,2] = "16_SANTA ANA"
$data[155,3] = 1
$data[155,4] = "LP1912"
$data[156,0] = "10:29:57"
$data[156,1] = "11:31"
$data[156,2] = "81_EL PELIGRO"
$data[156,3] = 62
$data[156,4] = "LP1912"
$data[157,0] = "09:35:26"
$data[157,1] = "11:32"
$data[157,2] = "81_EL PELIGRO"
$data[157,3] = 117
$data[157,4] = "LP1912"
$data[158,0] = "10:59:49"
$data[158,1] = "11:34"
$data[158,2] = "23_HERNANDEZ"
$data[158,3] = 35
$data[158,4] = "LP1912"
$data[159,0] = "10:29:57"
$data[159,1] = "11:35"
$data[159,2] = "11_ETCHEVERRY"
$data[159,3] = 66
$data[159,4] = "LP1912"
$data[160,0] = "10:29:57"
$data[160,1] = "11:40"
$data[160,2] = "10_OLMOS"
$data[160,3] = 71
$data[160,4] = "LP1912"
$data[161,0] = "10:29:57"
$data[161,1] = "11:41"
$data[161,2] = "17_ROMERO"
$data[161,3] = 72
$data[161,4] = "LP1912"
$data[162,0] = "10:59:49"
$data[162,1] = "11:42"
$data[162,2] = "11_ETCHEVERRY"
$data[162,3] = 43
$data[162,4] = "LP1912"
$data[163,0] = "10:59:49"
$data[163,1] = "11:43"
$data[163,2] = "10_OLMOS"
$data[163,3] = 44
$data[163,4] = "LP1912"
$data[164,0] = "11:30:45"
$data[164,1] = "11:44"
$data[164,2] = "11_ETCHEVERRY"
$data[164,3] = 14
$data[164,4] = "LP1912"
$data[165,0] = "10:29:57"
$data[165,1] = "11:50"
$data[165,2] = "215B_EL PATO"
$data[165,3] = 81
$data[165,4] = "LP1912"
$data[166,0] = "10:59:49"
$data[166,1] = "11:51"
$data[166,2] = "215B_EL PATO"
$data[166,3] = 52
$data[166,4] = "LP1912"
$data[167,0] = "10:59:49"
$data[167,1] = "11:52"
$data[167,2] = "15_ABASTO"
$data[167,3] = 53
$data[167,4] = "LP1912"
$data[168,0] = "11:56:55"
$data[168,1] = "11:56"
$data[168,2] = "16_SANTA ANA"
$data[168,3] = 0
$data[168,4] = "LP1912"
$data[169,0] = "10:29:57"
$data[169,1] = "11:58"
$data[169,2] = "225_GOMEZ"
$data[169,3] = 89
$data[169,4] = "LP1912"
$data[170,0] = "11:30:45"
$data[170,1] = "11:59"
$data[170,2] = "225_GOMEZ"
$data[170,3] = 29
$data[170,4] = "LP1912"
$data[171,0] = "10:29:57"
$data[171,1] = "12:01"
$data[171,2] = "84_COLONIA URQUIZA-ESC 49"
$data[171,3] = 92
$data[171,4] = "LP1912"
$data[172,0] = "10:59:49"
$data[172,1] = "12:02"
$data[172,2] = "84_COLONIA URQUIZA-ESC 49"
$data[172,3] = 63
$data[172,4] = "LP1912"
$data[173,0] = "11:30:45"
$data[173,1] = "12:04"
$data[173,2] = "23_HERNANDEZ"
$data[173,3] = 34
$data[173,4] = "LP1912"
$data[174,0] = "10:29:57"
$data[174,1] = "12:06"
$data[174,2] = "16_P MOR-SANTA ANA"
$data[174,3] = 97
$data[174,4] = "LP1912"
$data[175,0] = "10:59:49"
$data[175,1] = "12:06"
$data[175,2] = "14_ABASTO"
$data[175,3] = 67
$data[175,4] = "LP1912"
$data[176,0] = "10:59:49"
$data[176,1] = "12:10"
$data[176,2] = "10_OLMOS"
$data[176,3] = 71
$data[176,4] = "LP1912"
$data[177,0] = "11:56:55"
$data[177,1] = "12:12"
$data[177,2] = "10_OLMOS"
$data[177,3] = 16
$data[177,4] = "LP1912"
$data[178,0] = "10:29:57"
$data[178,1] = "12:13"
$data[178,2] = "17_ROMERO"
$data[178,3] = 104
$data[178,4] = "LP1912"
$data[179,0] = "11:56:55"
$data[179,1] = "12:14"
$data[179,2] = "17_ROMERO"
$data[179,3] = 18
$data[179,4] = "LP1912"
$data[180,0] = "10:29:57"
$data[180,1] = "12:15"
$data[180,2] = "14_ABASTO"
$data[180,3] = 106
$data[180,4] = "LP1912"
$data[181,0] = "10:29:57"
$data[181,1] = "12:20"
$data[181,2] = "215A_EL PATO"
$data[181,3] = 111
$data[181,4] = "LP1912"
$data[182,0] = "10:29:57"
$data[182,1] = "12:20"
$data[182,2] = "26_HERNANDEZ"
$data[182,3] = 111
$data[182,4] = "LP1912"
$data[183,0] = "10:59:49"
$data[183,1] = "12:20"
$data[183,2] = "14_ABASTO"
$data[183,3] = 81
$data[183,4] = "LP1912"
$data[184,0] = "10:59:49"
$data[184,1] = "12:21"
$data[184,2] = "26_HERNANDEZ"
$data[184,3] = 82
$data[184,4] = "LP1912"
$data[185,0] = "12:21:08"
$data[185,1] = "12:21"
$data[185,2] = "16_SANTA ANA"
$data[185,3] = 0
$data[185,4] = "LP1912"
$data[186,0] = "12:21:08"
$data[186,1] = "12:21"
$data[186,2] = "215A_EL PATO"
$data[186,3] = 0
$data[186,4] = "LP1912"
$data[187,0] = "10:59:49"
$data[187,1] = "12:30"
$data[187,2] = "17_ROMERO"
$data[187,3] = 91
$data[187,4] = "LP1912"
$data[188,0] = "11:56:55"
$data[188,1] = "12:34"
$data[188,2] = "23_HERNANDEZ"
$data[188,3] = 38
$data[188,4] = "LP1912"
$data[189,0] = "11:56:55"
$data[189,1] = "12:34"
$data[189,2] = "11_ETCHEVERRY"
$data[189,3] = 38
$data[189,4] = "LP1912"
$data[190,0] = "12:21:08"
$data[190,1] = "12:35"
$data[190,2] = "23_HERNANDEZ"
$data[190,3] = 14
$data[190,4] = "LP1912"
$data[191,0] = "12:21:08"
$data[191,1] = "12:35"
$data[191,2] = "11_ETCHEVERRY"
$data[191,3] = 14
$data[191,4] = "LP1912"
$data[192,0] = "10:59:49"
$data[192,1] = "12:36"
$data[192,2] = "27_EL RETIRO"
$data[192,3] = 97
$data[192,4] = "LP1912"
$data[193,0] = "12:21:08"
$data[193,1] = "12:37"
$data[193,2] = "27_EL RETIRO"
$data[193,3] = 16
$data[193,4] = "LP1912"
$data[194,0] = "10:59:49"
$data[194,1] = "12:38"
$data[194,2] = "17_179 Y 38"
$data[194,3] = 99
$data[194,4] = "LP1912"
$data[195,0] = "11:56:55"
$data[195,1] = "12:40"
$data[195,2] = "10_OLMOS"
$data[195,3] = 44
$data[195,4] = "LP1912"
$data[196,0] = "11:30:45"
$data[196,1] = "12:41"
$data[196,2] = "10_OLMOS"
$data[196,3] = 71
$data[196,4] = "LP1912"
$data[197,0] = "12:47:27"
$data[197,1] = "12:47"
$data[197,2] = "16_SANTA ANA"
$data[197,3] = 0
$data[197,4] = "LP1912"
$data[198,0] = "10:59:49"
$data[198,1] = "12:48"
$data[198,2] = "11_ETCHEVERRY"
$data[198,3] = 109
$data[198,4] = "LP1912"
$data[199,0] = "12:47:27"
$data[199,1] = "12:48"
$data[199,2] = "16_SANTA ANA"
$data[199,3] = 1
$data[199,4] = "LP1912"
$data[200,0] = "12:21:08"
$data[200,1] = "12:49"
$data[200,2] = "11_ETCHEVERRY"
$data[200,3] = 28
$data[200,4] = "LP1912"
$data[201,0] = "12:21:08"
$data[201,1] = "12:55"
$data[201,2] = "10_OLMOS"
$data[201,3] = 34
$data[201,4] = "LP1912"
$data[202,0] = "12:59:47"
$data[202,1] = "13:00"
$data[202,2] = "16_SANTA ANA"
$data[202,3] = 1
$data[202,4] = "LP1912"
$data[203,0] = "11:30:45"
$data[203,1] = "13:01"
$data[203,2] = "17_ROMERO"
$data[203,3] = 91
$data[203,4] = "LP1912"
$data[204,0] = "12:47:27"
$data[204,1] = "13:02"
$data[204,2] = "15_ABASTO"
$data[204,3] = 15
$data[204,4] = "LP1912"
$data[205,0] = "12:21:08"
$data[205,1] = "13:03"
$data[205,2] = "14_ABASTO"
$data[205,3] = 42
$data[205,4] = "LP1912"
$data[206,0] = "12:47:27"
$data[206,1] = "13:04"
$data[206,2] = "23_HERNANDEZ"
$data[206,3] = 17
$data[206,4] = "LP1912"
$data[207,0] = "12:59:47"
$data[207,1] = "13:05"
$data[207,2] = "23_HERNANDEZ"
$data[207,3] = 6
$data[207,4] = "LP1912"
$data[208,0] = "11:30:45"
$data[208,1] = "13:06"
$data[208,2] = "16_P MOR-SANTA ANA"
$data[208,3] = 96
$data[208,4] = "LP1912"
$data[209,0] = "12:21:08"
$data[209,1] = "13:07"
$data[209,2] = "16_P MOR-SANTA ANA"
$data[209,3] = 46
$data[209,4] = "LP1912"
$data[210,0] = "11:30:45"
$data[210,1] = "13:07"
$data[210,2] = "10_OLMOS"
$data[210,3] = 97
$data[210,4] = "LP1912"
$data[211,0] = "12:21:08"
$data[211,1] = "13:08"
$data[211,2] = "10_OLMOS"
$data[211,3] = 47
$data[211,4] = "LP1912"
$data[212,0] = "11:30:45"
$data[212,1] = "13:13"
$data[212,2] = "215D_EL PATO"
$data[212,3] = 103
$data[212,4] = "LP1912"
$data[213,0] = "12:21:08"
$data[213,1] = "13:14"
$data[213,2] = "215D_EL PATO"
$data[213,3] = 53
$data[213,4] = "LP1912"
$data[214,0] = "12:47:27"
$data[214,1] = "13:14"
$data[214,2] = "11_ETCHEVERRY"
$data[214,3] = 27
$data[214,4] = "LP1912"
$data[215,0] = "11:56:55"
$data[215,1] = "13:20"
$data[215,2] = "26_HERNANDEZ"
$data[215,3] = 84
$data[215,4] = "LP1912"
$data[216,0] = "11:30:45"
$data[216,1] = "13:21"
$data[216,2] = "26_HERNANDEZ"
$data[216,3] = 111
$data[216,4] = "LP1912"
$data[217,0] = "11:30:45"
$data[217,1] = "13:25"
$data[217,2] = "10_OLMOS"
$data[217,3] = 115
$data[217,4] = "LP1912"
$data[218,0] = "11:30:45"
$data[218,1] = "13:26"
$data[218,2] = "15_ABASTO"
$data[218,3] = 116
$data[218,4] = "LP1912"
$data[219,0] = "11:30:45"
$data[219,1] = "13:26"
$data[219,2] = "14_ABASTO"
$data[219,3] = 116
$data[219,4] = "LP1912"
$data[220,0] = "11:56:55"
$data[220,1] = "13:27"
$data[220,2] = "10_OLMOS"
$data[220,3] = 91
$data[220,4] = "LP1912"
$data[221,0] = "12:21:08"
$data[221,1] = "13:27"
$data[221,2] = "14_ABASTO"
$data[221,3] = 66
$data[221,4] = "LP1912"
$data[222,0] = "12:21:08"
$data[222,1] = "13:28"
$data[222,2] = "10_OLMOS"
$data[222,3] = 67
$data[222,4] = "LP1912"
$data[223,0] = "12:47:27"
$data[223,1] = "13:31"
$data[223,2] = "10_OLMOS"
$data[223,3] = 44
$data[223,4] = "LP1912"
$data[224,0] = "12:47:27"
$data[224,1] = "13:32"
$data[224,2] = "10_OLMOS"
$data[224,3] = 45
$data[224,4] = "LP1912"
$data[225,0] = "13:33:42"
$data[225,1] = "13:33"
$data[225,2] = "16_SANTA ANA"
$data[225,3] = 0
$data[225,4] = "LP1912"
$data[226,0] = "12:59:47"
$data[226,1] = "13:33"
$data[226,2] = "10_OLMOS"
$data[226,3] = 34
$data[226,4] = "LP1912"
$data[227,0] = "13:33:42"
$data[227,1] = "13:34"
$data[227,2] = "16_SANTA ANA"
$data[227,3] = 1
$data[227,4] = "LP1912"
$data[228,0] = "13:33:42"
$data[228,1] = "13:34"
$data[228,2] = "23_HERNANDEZ"
$data[228,3] = 1
$data[228,4] = "LP1912"
$data[229,0] = "11:56:55"
$data[229,1] = "13:36"
$data[229,2] = "15_ABASTO"
$data[229,3] = 100
$data[229,4] = "LP1912"
$data[230,0] = "13:33:42"
$data[230,1] = "13:38"
$data[230,2] = "14_ABASTO"
$data[230,3] = 5
$data[230,4] = "LP1912"
$data[231,0] = "11:56:55"
$data[231,1] = "13:46"
$data[231,2] = "17_ROMERO"
$data[231,3] = 110
$data[231,4] = "LP1912"
$data[232,0] = "12:59:47"
$data[232,1] = "13:50"
$data[232,2] = "11_ETCHEVERRY"
$data[232,3] = 51
$data[232,4] = "LP1912"
$data[233,0] = "11:56:55"
$data[233,1] = "13:50"
$data[233,2] = "215A_EL PATO"
$data[233,3] = 114
$data[233,4] = "LP1912"
$data[234,0] = "12:21:08"
$data[234,1] = "13:51"
$data[234,2] = "215A_EL PATO"
$data[234,3] = 90
$data[234,4] = "LP1912"
$data[235,0] = "11:56:55"
$data[235,1] = "13:55"
$data[235,2] = "225_GOMEZ"
$data[235,3] = 119
$data[235,4] = "LP1912"
$data[236,0] = "12:21:08"
$data[236,1] = "13:56"
$data[236,2] = "225_GOMEZ"
$data[236,3] = 95
$data[236,4] = "LP1912"
$data[237,0] = "12:59:47"
$data[237,1] = "13:56"
$data[237,2] = "16_P MOR-167 Y 521"
$data[237,3] = 57
$data[237,4] = "LP1912"
$data[238,0] = "12:47:27"
$data[238,1] = "13:58"
$data[238,2] = "16_P MOR-167 Y 521"
$data[238,3] = 71
$data[238,4] = "LP1912"
$data[239,0] = "13:59:06"
$data[239,1] = "13:59"
$data[239,2] = "16_SANTA ANA"
$data[239,3] = 0
$data[239,4] = "LP1912"
$data[240,0] = "13:59:06"
$data[240,1] = "14:00"
$data[240,2] = "16_SANTA ANA"
$data[240,3] = 1
$data[240,4] = "LP1912"
$data[241,0] = "12:21:08"
$data[241,1] = "14:00"
$data[241,2] = "16_P MOR-167 Y 521"
$data[241,3] = 99
$data[241,4] = "LP1912"
$data[242,0] = "12:21:08"
$data[242,1] = "14:04"
$data[242,2] = "17_ROMERO"
$data[242,3] = 103
$data[242,4] = "LP1912"
$data[243,0] = "13:33:42"
$data[243,1] = "14:04"
$data[243,2] = "23_HERNANDEZ"
$data[243,3] = 31
$data[243,4] = "LP1912"
$data[244,0] = "13:59:06"
$data[244,1] = "14:05"
$data[244,2] = "23_HERNANDEZ"
$data[244,3] = 6
$data[244,4] = "LP1912"
$data[245,0] = "12:21:08"
$data[245,1] = "14:08"
$data[245,2] = "23_HERNANDEZ"
$data[245,3] = 107
$data[245,4] = "LP1912"
$data[246,0] = "12:59:47"
$data[246,1] = "14:11"
$data[246,2] = "23_HERNANDEZ"
$data[246,3] = 72
$data[246,4] = "LP1912"
$data[247,0] = "13:33:42"
$data[247,1] = "14:12"
$data[247,2] = "15_ABASTO"
$data[247,3] = 39
$data[247,4] = "LP1912"
$data[248,0] = "12:47:27"
$data[248,1] = "14:16"
$data[248,2] = "27_EL RETIRO"
$data[248,3] = 89
$data[248,4] = "LP1912"
$data[249,0] = "12:21:08"
$data[249,1] = "14:17"
$data[249,2] = "27_EL RETIRO"
$data[249,3] = 116
$data[249,4] = "LP1912"
$data[250,0] = "12:59:47"
$data[250,1] = "14:19"
$data[250,2] = "215C_EL PATO"
$data[250,3] = 80
$data[250,4] = "LP1912"
$data[251,0] = "12:21:08"
$data[251,1] = "14:20"
$data[251,2] = "215C_EL PATO"
$data[251,3] = 119
$data[251,4] = "LP1912"
$data[252,0] = "12:47:27"
$data[252,1] = "14:21"
$data[252,2] = "26_HERNANDEZ"
$data[252,3] = 94
$data[252,4] = "LP1912"
$data[253,0] = "14:24:16"
$data[253,1] = "14:25"
$data[253,2] = "16_SANTA ANA"
$data[253,3] = 1
$data[253,4] = "LP1912"
$data[254,0] = "13:59:06"
$data[254,1] = "14:28"
$data[254,2] = "15_ABASTO"
$data[254,3] = 29
$data[254,4] = "LP1912"
$data[255,0] = "14:24:16"
$data[255,1] = "14:35"
$data[255,2] = "23_HERNANDEZ"
$data[255,3] = 11
$data[255,4] = "LP1912"
$data[256,0] = "14:24:16"
$data[256,1] = "14:44"
$data[256,2] = "15_ABASTO"
$data[256,3] = 20
$data[256,4] = "LP1912"
$data[257,0] = "13:33:42"
$data[257,1] = "14:44"
$data[257,2] = "14_ABASTO"
$data[257,3] = 71
$data[257,4] = "LP1912"
$data[258,0] = "12:47:27"
$data[258,1] = "14:45"
$data[258,2] = "14_ABASTO"
$data[258,3] = 118
$data[258,4] = "LP1912"
$data[259,0] = "14:45:17"
$data[259,1] = "14:45"
$data[259,2] = "15_ABASTO"
$data[259,3] = 0
$data[259,4] = "LP1912"
$data[260,0] = "14:45:17"
$data[260,1] = "14:46"
$data[260,2] = "16_SANTA ANA"
$data[260,3] = 1
$data[260,4] = "LP1912"
$data[261,0] = "14:56:20"
$data[261,1] = "14:56"
$data[261,2] = "16_SANTA ANA"
$data[261,3] = 0
$data[261,4] = "LP1912"
$data[262,0] = "12:59:47"
$data[262,1] = "14:56"
$data[262,2] = "16_P MOR-SANTA ANA"
$data[262,3] = 117
$data[262,4] = "LP1912"
$data[263,0] = "13:59:06"
$data[263,1] = "14:57"
$data[263,2] = "16_P MOR-SANTA ANA"
$data[263,3] = 58
$data[263,4] = "LP1912"
$data[264,0] = "12:59:47"
$data[264,1] = "14:58"
$data[264,2] = "215B_EL PATO"
$data[264,3] = 119
$data[264,4] = "LP1912"
$data[265,0] = "13:33:42"
$data[265,1] = "15:00"
$data[265,2] = "81_EL PELIGRO"
$data[265,3] = 87
$data[265,4] = "LP1912"
$data[266,0] = "13:33:42"
$data[266,1] = "15:05"
$data[266,2] = "10_OLMOS"
$data[266,3] = 92
$data[266,4] = "LP1912"
$data[267,0] = "14:45:17"
$data[267,1] = "15:05"
$data[267,2] = "23_HERNANDEZ"
$data[267,3] = 20
$data[267,4] = "LP1912"
$data[268,0] = "13:59:06"
$data[268,1] = "15:10"
$data[268,2] = "17_ROMERO"
$data[268,3] = 71
$data[268,4] = "LP1912"
$data[269,0] = "13:33:42"
$data[269,1] = "15:13"
$data[269,2] = "11_ETCHEVERRY"
$data[269,3] = 100
$data[269,4] = "LP1912"
$data[270,0] = "13:59:06"
$data[270,1] = "15:14"
$data[270,2] = "11_ETCHEVERRY"
$data[270,3] = 75
$data[270,4] = "LP1912"
$data[271,0] = "13:33:42"
$data[271,1] = "15:17"
$data[271,2] = "26_HERNANDEZ"
$data[271,3] = 104
$data[271,4] = "LP1912"
$data[272,0] = "14:56:20"
$data[272,1] = "15:17"
$data[272,2] = "16_SANTA ANA"
$data[272,3] = 21
$data[272,4] = "LP1912"
$data[273,0] = "13:59:06"
$data[273,1] = "15:18"
$data[273,2] = "26_HERNANDEZ"
$data[273,3] = 79
$data[273,4] = "LP1912"
$data[274,0] = "14:56:20"
$data[274,1] = "15:20"
$data[274,2] = "15_ABASTO"
$data[274,3] = 24
$data[274,4] = "LP1912"
$data[275,0] = "14:24:16"
$data[275,1] = "15:21"
$data[275,2] = "26_HERNANDEZ"
$data[275,3] = 57
$data[275,4] = "LP1912"
$data[276,0] = "15:22:17"
$data[276,1] = "15:22"
$data[276,2] = "16_SANTA ANA"
$data[276,3] = 0
$data[276,4] = "LP1912"
$data[277,0] = "15:22:17"
$data[277,1] = "15:22"
$data[277,2] = "26_HERNANDEZ"
$data[277,3] = 0
$data[277,4] = "LP1912"
$data[278,0] = "14:24:16"
$data[278,1] = "15:32"
$data[278,2] = "84_COLONIA URQUIZA-ESC 49"
$data[278,3] = 68
$data[278,4] = "LP1912"
$data[279,0] = "13:59:06"
$data[279,1] = "15:35"
$data[279,2] = "23_HERNANDEZ"
$data[279,3] = 96
$data[279,4] = "LP1912"
$data[280,0] = "13:59:06"
$data[280,1] = "15:37"
$data[280,2] = "10_OLMOS"
$data[280,3] = 98
$data[280,4] = "LP1912"
$data[281,0] = "14:24:16"
$data[281,1] = "15:38"
$data[281,2] = "23_HERNANDEZ"
$data[281,3] = 74
$data[281,4] = "LP1912"
$data[282,0] = "14:45:17"
$data[282,1] = "15:38"
$data[282,2] = "215A_EL PATO"
$data[282,3] = 53
$data[282,4] = "LP1912"
$data[283,0] = "14:56:20"
$data[283,1] = "15:38"
$data[283,2] = "10_OLMOS"
$data[283,3] = 42
$data[283,4] = "LP1912"
$data[284,0] = "13:59:06"
$data[284,1] = "15:39"
$data[284,2] = "215A_EL PATO"
$data[284,3] = 100
$data[284,4] = "LP1912"
$data[285,0] = "14:56:20"
$data[285,1] = "15:45"
$data[285,2] = "14_ABASTO"
$data[285,3] = 49
$data[285,4] = "LP1912"
$data[286,0] = "14:24:16"
$data[286,1] = "15:46"
$data[286,2] = "14_ABASTO"
$data[286,3] = 82
$data[286,4] = "LP1912"
$data[287,0] = "14:56:20"
$data[287,1] = "15:46"
$data[287,2] = "16_P MOR-167 Y 521"
$data[287,3] = 50
$data[287,4] = "LP1912"
$data[288,0] = "13:59:06"
$data[288,1] = "15:47"
$data[288,2] = "16_P MOR-167 Y 521"
$data[288,3] = 108
$data[288,4] = "LP1912"
$data[289,0] = "13:59:06"
$data[289,1] = "15:48"
$data[289,2] = "14_ABASTO"
$data[289,3] = 109
$data[289,4] = "LP1912"
$data[290,0] = "14:56:20"
$data[290,1] = "15:53"
$data[290,2] = "11_ETCHEVERRY"
$data[290,3] = 57
$data[290,4] = "LP1912"
$data[291,0] = "13:59:06"
$data[291,1] = "15:54"
$data[291,2] = "11_ETCHEVERRY"
$data[291,3] = 115
$data[291,4] = "LP1912"
$data[292,0] = "15:53:28"
$data[292,1] = "15:54"
$data[292,2] = "16_SANTA ANA"
$data[292,3] = 1
$data[292,4] = "LP1912"
$data[293,0] = "15:22:17"
$data[293,1] = "15:55"
$data[293,2] = "16_SANTA ANA"
$data[293,3] = 33
$data[293,4] = "LP1912"
$data[294,0] = "15:53:28"
$data[294,1] = "15:56"
$data[294,2] = "27_EL RETIRO"
$data[294,3] = 3
$data[294,4] = "LP1912"
$data[295,0] = "14:24:16"
$data[295,1] = "15:56"
$data[295,2] = "17_ROMERO"
$data[295,3] = 92
$data[295,4] = "LP1912"
$data[296,0] = "13:59:06"
$data[296,1] = "15:57"
$data[296,2] = "27_EL RETIRO"
$data[296,3] = 118
$data[296,4] = "LP1912"
$data[297,0] = "15:22:17"
$data[297,1] = "16:01"
$data[297,2] = "10_OLMOS"
$data[297,3] = 39
$data[297,4] = "LP1912"
$data[298,0] = "15:53:28"
$data[298,1] = "16:02"
$data[298,2] = "16_SANTA ANA"
$data[298,3] = 9
$data[298,4] = "LP1912"
$data[299,0] = "15:53:28"
$data[299,1] = "16:04"
$data[299,2] = "23_HERNANDEZ"
$data[299,3] = 11
$data[299,4] = "LP1912"
$data[300,0] = "15:22:17"
$data[300,1] = "16:05"
$data[300,2] = "23_HERNANDEZ"
$data[300,3] = 43
$data[300,4] = "LP1912"
$data[301,0] = "14:56:20"
$data[301,1] = "16:08"
$data[301,2] = "14_ABASTO"
$data[301,3] = 72
$data[301,4] = "LP1912"
$data[302,0] = "14:45:17"
$data[302,1] = "16:09"
$data[302,2] = "14_ABASTO"
$data[302,3] = 84
$data[302,4] = "LP1912"
$data[303,0] = "14:24:16"
$data[303,1] = "16:15"
$data[303,2] = "225_C ROCA-H SUR"
$data[303,3] = 111
$data[303,4] = "LP1912"
$data[304,0] = "16:13:37"
$data[304,1] = "16:15"
$data[304,2] = "16_SANTA ANA"
$data[304,3] = 2
$data[304,4] = "LP1912"
$data[305,0] = "16:13:37"
$data[305,1] = "16:19"
$data[305,2] = "215C_EL PATO"
$data[305,3] = 6
$data[305,4] = "LP1912"
$data[306,0] = "14:24:16"
$data[306,1] = "16:20"
$data[306,2] = "215C_EL PATO"
$data[306,3] = 116
$data[306,4] = "LP1912"
$data[307,0] = "14:24:16"
$data[307,1] = "16:21"
$data[307,2] = "26_HERNANDEZ"
$data[307,3] = 117
$data[307,4] = "LP1912"
$data[308,0] = "16:13:37"
$data[308,1] = "16:26"
$data[308,2] = "16_SANTA ANA"
$data[308,3] = 13
$data[308,4] = "LP1912"
$data[309,0] = "15:53:28"
$data[309,1] = "16:29"
$data[309,2] = "10_OLMOS"
$data[309,3] = 36
$data[309,4] = "LP1912"
$data[310,0] = "14:45:17"
$data[310,1] = "16:30"
$data[310,2] = "15_ABASTO"
$data[310,3] = 105
$data[310,4] = "LP1912"
$data[311,0] = "16:31:51"
$data[311,1] = "16:31"
$data[311,2] = "16_SANTA ANA"
$data[311,3] = 0
$data[311,4] = "LP1912"
$data[312,0] = "15:22:17"
$data[312,1] = "16:32"
$data[312,2] = "14_ABASTO"
$data[312,3] = 70
$data[312,4] = "LP1912"
$data[313,0] = "15:53:28"
$data[313,1] = "16:34"
$data[313,2] = "23_HERNANDEZ"
$data[313,3] = 41
$data[313,4] = "LP1912"
$data[314,0] = "15:53:28"
$data[314,1] = "16:36"
$data[314,2] = "11_ETCHEVERRY"
$data[314,3] = 43
$data[314,4] = "LP1912"
$data[315,0] = "15:22:17"
$data[315,1] = "16:37"
$data[315,2] = "11_ETCHEVERRY"
$data[315,3] = 75
$data[315,4] = "LP1912"
$data[316,0] = "15:22:17"
$data[316,1] = "16:40"
$data[316,2] = "17_ROMERO"
$data[316,3] = 78
$data[316,4] = "LP1912"
$data[317,0] = "16:31:51"
$data[317,1] = "16:42"
$data[317,2] = "225_GOMEZ"
$data[317,3] = 11
$data[317,4] = "LP1912"
$data[318,0] = "14:56:20"
$data[318,1] = "16:42"
$data[318,2] = "16_P MOR-SANTA ANA"
$data[318,3] = 106
$data[318,4] = "LP1912"
$data[319,0] = "14:45:17"
$data[319,1] = "16:43"
$data[319,2] = "16_P MOR-SANTA ANA"
$data[319,3] = 118
$data[319,4] = "LP1912"
$data[320,0] = "14:45:17"
$data[320,1] = "16:43"
$data[320,2] = "225_GOMEZ"
$data[320,3] = 118
$data[320,4] = "LP1912"
$data[321,0] = "16:45:31"
$data[321,1] = "16:45"
$data[321,2] = "16_SANTA ANA"
$data[321,3] = 0
$data[321,4] = "LP1912"
$data[322,0] = "15:22:17"
$data[322,1] = "16:48"
$data[322,2] = "15_ABASTO"
$data[322,3] = 86
$data[322,4] = "LP1912"
$data[323,0] = "15:53:28"
$data[323,1] = "16:50"
$data[323,2] = "14_ABASTO"
$data[323,3] = 57
$data[323,4] = "LP1912"
$data[324,0] = "16:52:36"
$data[324,1] = "16:52"
$data[324,2] = "10_OLMOS"
$data[324,3] = 0
$data[324,4] = "LP1912"
$data[325,0] = "16:52:36"
$data[325,1] = "16:53"
$data[325,2] = "16_SANTA ANA"
$data[325,3] = 1
$data[325,4] = "LP1912"
$data[326,0] = "16:31:51"
$data[326,1] = "16:56"
$data[326,2] = "10_OLMOS"
$data[326,3] = 25
$data[326,4] = "LP1912"
$data[327,0] = "15:22:17"
$data[327,1] = "16:56"
$data[327,2] = "17_179 Y 38"
$data[327,3] = 94
$data[327,4] = "LP1912"
$data[328,0] = "16:13:37"
$data[328,1] = "16:57"
$data[328,2] = "10_OLMOS"
$data[328,3] = 44
$data[328,4] = "LP1912"
$data[329,0] = "15:22:17"
$data[329,1] = "17:04"
$data[329,2] = "215A_EL PATO"
$data[329,3] = 102
$data[329,4] = "LP1912"
$data[330,0] = "16:13:37"
$data[330,1] = "17:04"
$data[330,2] = "11_ETCHEVERRY"
$data[330,3] = 51
$data[330,4] = "LP1912"
$data[331,0] = "16:13:37"
$data[331,1] = "17:04"
$data[331,2] = "23_HERNANDEZ"
$data[331,3] = 51
$data[331,4] = "LP1912"
$data[332,0] = "16:52:36"
$data[332,1] = "17:05"
$data[332,2] = "23_HERNANDEZ"
$data[332,3] = 13
$data[332,4] = "LP1912"
$data[333,0] = "16:45:31"
$data[333,1] = "17:06"
$data[333,2] = "23_HERNANDEZ"
$data[333,3] = 21
$data[333,4] = "LP1912"
$data[334,0] = "16:31:51"
$data[334,1] = "17:09"
$data[334,2] = "10_OLMOS"
$data[334,3] = 38
$data[334,4] = "LP1912"
$data[335,0] = "16:45:31"
$data[335,1] = "17:10"
$data[335,2] = "10_OLMOS"
$data[335,3] = 25
$data[335,4] = "LP1912"
$data[336,0] = "16:45:31"
$data[336,1] = "17:16"
$data[336,2] = "11_ETCHEVERRY"
$data[336,3] = 31
$data[336,4] = "LP1912"
$data[337,0] = "16:31:51"
$data[337,1] = "17:20"
$data[337,2] = "26_HERNANDEZ"
$data[337,3] = 49
$data[337,4] = "LP1912"
$data[338,0] = "16:31:51"
$data[338,1] = "17:20"
$data[338,2] = "16_SANTA ANA"
$data[338,3] = 49
$data[338,4] = "LP1912"
$data[339,0] = "15:53:28"
$data[339,1] = "17:21"
$data[339,2] = "26_HERNANDEZ"
$data[339,3] = 88
$data[339,4] = "LP1912"
$data[340,0] = "15:53:28"
$data[340,1] = "17:24"
$data[340,2] = "84_COLONIA URQUIZA-ESC 49"
$data[340,3] = 91
$data[340,4] = "LP1912"
$data[341,0] = "15:53:28"
$data[341,1] = "17:28"
$data[341,2] = "14_ABASTO"
$data[341,3] = 95
$data[341,4] = "LP1912"
$data[342,0] = "16:45:31"
$data[342,1] = "17:34"
$data[342,2] = "23_HERNANDEZ"
$data[342,3] = 49
$data[342,4] = "LP1912"
$data[343,0] = "15:53:28"
$data[343,1] = "17:36"
$data[343,2] = "27_EL RETIRO"
$data[343,3] = 103
$data[343,4] = "LP1912"
$data[344,0] = "15:53:28"
$data[344,1] = "17:38"
$data[344,2] = "17_ROMERO"
$data[344,3] = 105
$data[344,4] = "LP1912"
$data[345,0] = "15:53:28"
$data[345,1] = "17:40"
$data[345,2] = "215B_EL PATO"
$data[345,3] = 107
$data[345,4] = "LP1912"
$data[346,0] = "16:13:37"
$data[346,1] = "17:40"
$data[346,2] = "17_ROMERO"
$data[346,3] = 87
$data[346,4] = "LP1912"
$data[347,0] = "16:31:51"
$data[347,1] = "17:45"
$data[347,2] = "15_ABASTO"
$data[347,3] = 74
$data[347,4] = "LP1912"
$data[348,0] = "15:53:28"
$data[348,1] = "17:50"
$data[348,2] = "16_P MOR-167 Y 521"
$data[348,3] = 117
$data[348,4] = "LP1912"
$data[349,0] = "15:53:28"
$data[349,1] = "17:52"
$data[349,2] = "81_EL PELIGRO"
$data[349,3] = 119
$data[349,4] = "LP1912"
$data[350,0] = "16:13:37"
$data[350,1] = "18:04"
$data[350,2] = "17_ROMERO"
$data[350,3] = 111
$data[350,4] = "LP1912"
$data[351,0] = "16:52:36"
$data[351,1] = "18:08"
$data[351,2] = "14_ABASTO"
$data[351,3] = 76
$data[351,4] = "LP1912"
$data[352,0] = "16:31:51"
$data[352,1] = "18:20"
$data[352,2] = "26_HERNANDEZ"
$data[352,3] = 109
$data[352,4] = "LP1912"
$data[353,0] = "16:45:31"
$data[353,1] = "18:21"
$data[353,2] = "26_HERNANDEZ"
$data[353,3] = 96
$data[353,4] = "LP1912"
$data[354,0] = "16:31:51"
$data[354,1] = "18:27"
$data[354,2] = "215C_EL PATO"
$data[354,3] = 116
$data[354,4] = "LP1912"
$data[355,0] = "16:45:31"
$data[355,1] = "18:28"
$data[355,2] = "215C_EL PATO"
$data[355,3] = 103
$data[355,4] = "LP1912"
$data[356,0] = "16:45:31"
$data[356,1] = "18:32"
$data[356,2] = "11X44_ETCHEVERRY"
$data[356,3] = 107
$data[356,4] = "LP1912"
$data[357,0] = "16:45:31"
$data[357,1] = "18:40"
$data[357,2] = "15_ABASTO"
$data[357,3] = 115
$data[357,4] = "LP1912"
$data[358,0] = "16:52:36"
$data[358,1] = "18:48"
$data[358,2] = "14X44_ABASTO"
$data[358,3] = 116
$data[358,4] = "LP1912"
$ws.Range("A6:E364").Value = $data

